# "se esta realizando la diapositiva de clase 6"
# Mark column "6" (class 6, columns I/J) attendance as present ("p")
# for every student row on the attendance sheet, and move the
# viewport/selection to reflect where the teacher is currently working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Select()

# Row 24 was missing the "p" mark in column H (class 4) - fill it in too.
$ws.Range("H24").Value = "p"

# Mark columns I (class 5) and J (class 6) present for every student (rows 3-29).
For ($r = 3; $r -le 29; $r++) {
    $ws.Range("I$r").Value = "p"
    $ws.Range("J$r").Value = "p"
}

# Reflect the scroll position / active selection recorded when the class-6
# slide was being prepared.
$ws.Range("A8").Select()
$ws.Range("K11").Select()
